# kNN classifier with Mahalanobis distance is added.
#
# The "Control" sheet's sample table had two rows removed (ID=65 and
# ID=79), which shifts every following row up and shrinks the table
# from A1:C40 to A1:C38. The "Pathology" sheet's selection also moves,
# and "Control" becomes the active tab/sheet.

$wb = $excel.ActiveWorkbook

$pathology = $wb.Worksheets.Item("Pathology")
$control   = $wb.Worksheets.Item("Control")

# Remove the two data rows (ID=65 at row 26, ID=79 at row 30 in the
# original layout). Deleting the higher-numbered row first keeps the
# row 26 index valid for the second delete.
[void]$control.Rows.Item(30).Delete()
[void]$control.Rows.Item(26).Delete()

# Pathology ends up with its selection on C15 (no longer the active tab).
[void]$pathology.Activate()
[void]$pathology.Range("C15").Select()

# Control becomes the active tab, scrolled/selected at row 29 (entire row).
[void]$control.Activate()
[void]$control.Rows.Item(29).Select()
